$d = $word.ActiveDocument

# Change 1: "eslé tient longtemps" -> "eslé tient long temps"
$d.Content.Find.Execute("eslé tient longtemps", $true, $false, $false, $false, $false,
                         $true, 1, $false, "eslé tient long temps", 2)

# Change 2: " non pas pres du col," -> " non pas pr" + "é" (different formatting) + "s du col,"
$d.Content.Find.Execute(" non pas pres du col,", $true, $false, $false, $false, $false,
                         $true, 1, $false, " non pas prés du col,", 2)
